# The deck's date placeholders (on the slide master and every slide layout)
# hold an auto-updating "datetime1" field whose cached display text is
# "12/7/2022". PowerPoint recomputed/re-cached that text to "12/12/2022"
# (e.g. on a later open/save) without any other visible change, so we just
# need to refresh the cached text of every Date Placeholder shape.

$p = $ppt.ActivePresentation
$newDate = "12/12/2022"
$ppPlaceholderDate = 16

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master's own Date Placeholder.
Update-DateShapes $p.SlideMaster.Shapes

# Every slide layout's Date Placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DateShapes $layouts.Item($L).Shapes
}
